$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H58").Value = 78944.25
$ws.Range("I58").Value = 169519.67
$ws.Range("J58").Value = 24599
$ws.Range("K58").Value = 508559.01
$ws.Range("L58").Value = 73797
$ws.Range("M58").Value = -508409.01
$ws.Range("N58").Value = -74097

$ws.Range("H96").Value = 1127.2858
$ws.Range("I96").Value = 982.6667
$ws.Range("J96").Value = 1995
$ws.Range("K96").Value = 2948.0001
$ws.Range("L96").Value = 5985
$ws.Range("M96").Value = -1575.0001
$ws.Range("N96").Value = -8731

$ws.Range("H116").Value = 62508450
$ws.Range("I116").Value = 83337940
$ws.Range("J116").Value = 20000
$ws.Range("K116").Value = 83337940
$ws.Range("L116").Value = 20000
$ws.Range("M116").Value = -83334498
$ws.Range("N116").Value = -26884

$ws.Range("H138").Value = 1856.2858
$ws.Range("I138").Value = 1410.7059
$ws.Range("J138").Value = 3750
$ws.Range("K138").Value = 4232.1177
$ws.Range("L138").Value = 11250
$ws.Range("M138").Value = 907.8823000000002
$ws.Range("N138").Value = -21530

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 4194.4287
$ws.Range("I2").Value = 1246.7142
$ws.Range("J2").Value = 7142.143
$ws.Range("K2").Value = 1246.7142
$ws.Range("L2").Value = 7142.143
$ws.Range("M2").Value = -1133.7142
$ws.Range("N2").Value = -7368.143

$ws.Range("H32").Value = 1474088.4
$ws.Range("I32").Value = 1565938.2
$ws.Range("K32").Value = 1565938.2
$ws.Range("M32").Value = -1565651.2

$ws.Range("H64").Value = 0
$ws.Range("J64").Value = 0
$ws.Range("L64").Value = 0
$ws.Range("N64").ClearContents()

$ws.Range("H67").Value = 0
$ws.Range("J67").Value = 0
$ws.Range("L67").Value = 0
$ws.Range("N67").ClearContents()

$ws.Range("H74").Value = 48140.39
$ws.Range("I74").Value = 82260.75
$ws.Range("K74").Value = 82260.75
$ws.Range("M74").Value = -81386.75

$ws.Range("H77").Value = 48140.39
$ws.Range("I77").Value = 82260.75
$ws.Range("K77").Value = 411303.75
$ws.Range("M77").Value = -406935.75

$ws.Range("H116").Value = 4194.4287
$ws.Range("I116").Value = 1246.7142
$ws.Range("J116").Value = 7142.143
$ws.Range("K116").Value = 1246.7142
$ws.Range("L116").Value = 7142.143
$ws.Range("M116").Value = 1047.2858
$ws.Range("N116").Value = -11730.143

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 4194.4287
$ws.Range("I3").Value = 1246.7142
$ws.Range("J3").Value = 7142.143
$ws.Range("K3").Value = 1246.7142
$ws.Range("L3").Value = 7142.143
$ws.Range("M3").Value = -1132.7142
$ws.Range("N3").Value = -7370.143

$ws.Range("H62").Value = 43000
$ws.Range("J62").Value = 43000
$ws.Range("L62").Value = 43000
$ws.Range("N62").Value = -44372

$ws.Range("H65").Value = 43000
$ws.Range("J65").Value = 43000
$ws.Range("L65").Value = 129000
$ws.Range("N65").Value = -135864

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 54.77778
$ws.Range("I7").Value = 31
$ws.Range("J7").Value = 73.8
$ws.Range("K7").Value = 31
$ws.Range("L7").Value = 73.8
$ws.Range("M7").Value = 82
$ws.Range("N7").Value = -299.8

$ws.Range("H16").Value = 2580.5588
$ws.Range("I16").Value = 1209.1818
$ws.Range("J16").Value = 5094.75
$ws.Range("K16").Value = 1209.1818
$ws.Range("L16").Value = 5094.75
$ws.Range("M16").Value = -922.1818000000001
$ws.Range("N16").Value = -5668.75

$ws.Range("H31").Value = 9017042
$ws.Range("I31").Value = 3398.0588
$ws.Range("K31").Value = 3398.0588
$ws.Range("M31").Value = -3103.0588

$ws.Range("H34").Value = 9017042
$ws.Range("I34").Value = 3398.0588
$ws.Range("K34").Value = 3398.0588
$ws.Range("M34").Value = -3196.0588

$ws.Range("H113").Value = 2580.5588
$ws.Range("I113").Value = 1209.1818
$ws.Range("J113").Value = 5094.75
$ws.Range("K113").Value = 1209.1818
$ws.Range("L113").Value = 5094.75
$ws.Range("M113").Value = 960.8181999999999
$ws.Range("N113").Value = -9434.75

$ws.Range("H134").Value = 5234.0933
$ws.Range("I134").Value = 1470.3636
$ws.Range("K134").Value = 4411.0908
$ws.Range("M134").Value = -1876.0908

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H47").Value = 442.22223
$ws.Range("I47").Value = 410
$ws.Range("K47").Value = 1230
$ws.Range("M47").Value = -799

$ws.Range("H48").Value = 9200
$ws.Range("J48").Value = 9200
$ws.Range("L48").Value = 27600
$ws.Range("N48").Value = -28100

$ws.Range("H49").Value = 0
$ws.Range("I49").Value = 0
$ws.Range("K49").Value = 0
$ws.Range("M49").ClearContents()

$ws.Range("H82").Value = 8166.6665
$ws.Range("I82").Value = 2250
$ws.Range("J82").Value = 20000
$ws.Range("K82").Value = 6750
$ws.Range("L82").Value = 60000
$ws.Range("M82").Value = -6344
$ws.Range("N82").Value = -60812

$ws.Range("H85").Value = 8166.6665
$ws.Range("I85").Value = 2250
$ws.Range("J85").Value = 20000
$ws.Range("K85").Value = 6750
$ws.Range("L85").Value = 60000
$ws.Range("M85").Value = -5346
$ws.Range("N85").Value = -62808

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H53").Value = 5000
$ws.Range("J53").Value = 5000
$ws.Range("L53").Value = 5000
$ws.Range("N53").Value = -6262

$ws.Range("H102").Value = 3137.0715
$ws.Range("I102").Value = 2922.3635
$ws.Range("K102").Value = 2922.3635
$ws.Range("M102").Value = -1300.3635

$ws.Range("H107").Value = 727854.4399999999
$ws.Range("I107").Value = 1143128.6
$ws.Range("J107").Value = 1124.75
$ws.Range("K107").Value = 1143128.6
$ws.Range("L107").Value = 1124.75
$ws.Range("M107").Value = -1141208.6
$ws.Range("N107").Value = -4964.75

$ws.Range("H122").Value = 2339149.8
$ws.Range("I122").Value = 3450961.5
$ws.Range("J122").Value = 4345.6
$ws.Range("K122").Value = 10352884.5
$ws.Range("L122").Value = 13036.8
$ws.Range("M122").Value = -10350434.5
$ws.Range("N122").Value = -17936.8

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 3278
$ws.Range("I22").Value = 500.5
$ws.Range("K22").Value = 500.5
$ws.Range("M22").Value = -205.5

$ws.Range("H27").Value = 3278
$ws.Range("I27").Value = 500.5
$ws.Range("K27").Value = 500.5
$ws.Range("M27").Value = -393.5

$ws.Range("H29").Value = 0
$ws.Range("I29").Value = 0
$ws.Range("K29").Value = 0
$ws.Range("M29").ClearContents()

$ws.Range("H46").Value = 2755.8
$ws.Range("I46").Value = 500.9
$ws.Range("J46").Value = 3657.76
$ws.Range("K46").Value = 500.9
$ws.Range("L46").Value = 3657.76
$ws.Range("M46").Value = -312.9
$ws.Range("N46").Value = -4033.76

$ws.Range("H55").Value = 34483120
$ws.Range("J55").Value = 594.0625
$ws.Range("L55").Value = 594.0625
$ws.Range("N55").Value = -940.0625

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 951.8182
$ws.Range("I100").Value = 764.8461
$ws.Range("J100").Value = 1221.8889
$ws.Range("K100").Value = 1529.6922
$ws.Range("L100").Value = 2443.7778
$ws.Range("M100").Value = -988.6922
$ws.Range("N100").Value = -3525.7778

